$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the instruction text in column B for rows 2-9 (row 1 is the header, unchanged)
$ws.Range("B2").Value = "Position the U-Channel bracket part number NYK:Z610399470001 in the housing "
$ws.Range("B3").Value = "Insert the two ¼-20 x 1.0” bolts. Tighten the bolts with a 3/8” socket or nut driver and torque to 50 in-lbs"
$ws.Range("B4").Value = "Install the manifold assembly part number NYK:9000903920001 in the U-Channel bracket as shown"
$ws.Range("B5").Value = "Insert three of the NYK:916001420000 ¼-20 x ½” bolts and tighten with a 3/8” socket or nut driver and torque each bolt to 50 in-lbs"
$ws.Range("B6").Value = "Take the battery part number NYK:Z801011870000 and install in the housing with the protruding side down . "
$ws.Range("B7").Value = "Take battery cushion NYK:Z916001440000 and apply it to the battery cover."
$ws.Range("B8").Value = "Place the Battery cover NYK:Z640399550001 in assembly."
$ws.Range("B9").Value = "Install the four M-F standoffs. Tighten with a 3/8” socket or nut driver and torque to 30 in-lbs."

# Set the row heights to match the re-wrapped text (rows 2-9) before removing the
# now-unused trailing rows, since heights stay tied to row position on delete.
$ws.Range("A2:C2").RowHeight = 45
$ws.Range("A3:C3").RowHeight = 45
$ws.Range("A4:C4").RowHeight = 45
$ws.Range("A5:C5").RowHeight = 60
$ws.Range("A6:C6").RowHeight = 45
$ws.Range("A7:C7").RowHeight = 45
$ws.Range("A8:C8").RowHeight = 30
$ws.Range("A9:C9").RowHeight = 45

# Delete the now-unused rows 10-18 (content was condensed into rows 2-9)
$ws.Range("A10:C18").EntireRow.Delete()

# Match the final selection/view state
$ws.Range("A9").Select() | Out-Null
